$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D6").Value = -7.855099999999993
$ws.Range("C7").Value = -12.59850000000001
$ws.Range("E7").Value = 15.93949999999999
$ws.Range("A9").Value = -21.77850000000001
$ws.Range("E10").Value = 16.66739999999999
$ws.Range("C12").Value = -10.6413
$ws.Range("A13").Value = -22.2079
$ws.Range("E13").Value = 16.7347
$ws.Range("C14").Value = -14.57129999999999
$ws.Range("D15").Value = -8.891999999999996
$ws.Range("A16").Value = -21.9118
$ws.Range("E16").Value = 16.09239999999999
$ws.Range("A18").Value = -22.34560000000001
$ws.Range("C19").Value = -12.26640000000001
$ws.Range("A20").Value = -21.54219999999997
$ws.Range("E20").Value = 15.16579999999999
$ws.Range("E24").Value = 16.75200000000001
$ws.Range("A26").Value = -21.70529999999999
$ws.Range("C26").Value = -13.11860000000001
$ws.Range("A27").Value = -21.81509999999999
$ws.Range("C27").Value = -12.7505
$ws.Range("D28").Value = -8.608899999999997
$ws.Range("A29").Value = -20.42549999999997
$ws.Range("C29").Value = -11.4724
$ws.Range("E32").Value = 16.90139999999999
$ws.Range("D33").Value = -7.8006
$ws.Range("A35").Value = -21.0057
$ws.Range("D35").Value = -8.440199999999994
$ws.Range("A36").Value = -21.2445
$ws.Range("C37").Value = -13.6908
$ws.Range("C38").Value = -13.08740000000001
$ws.Range("D38").Value = -9.229299999999993
$ws.Range("E39").Value = 16.4564
$ws.Range("D43").Value = -8.616499999999998
$ws.Range("D44").Value = -7.905600000000001
$ws.Range("A45").Value = -21.6824
$ws.Range("D45").Value = -7.443499999999999
$ws.Range("C47").Value = -11.4771
$ws.Range("D47").Value = -8.173399999999999
$ws.Range("E47").Value = 16.64190000000001
$ws.Range("E48").Value = 17.182
$ws.Range("C51").Value = -12.9533
$ws.Range("D51").Value = -7.697199999999999
$ws.Range("C52").Value = -10.99120000000001
$ws.Range("E52").Value = 17.20060000000001
$ws.Range("D54").Value = -8.116899999999994
$ws.Range("A55").Value = -22.17530000000001
$ws.Range("C55").Value = -13.01549999999999
$ws.Range("E56").Value = 16.5496
$ws.Range("A57").Value = -22.0744
$ws.Range("D57").Value = -8.146599999999998
$ws.Range("D62").Value = -8.293199999999999
$ws.Range("D63").Value = -7.773300000000002
$ws.Range("D67").Value = -5.865999999999998
$ws.Range("A69").Value = -21.77140000000001
$ws.Range("C69").Value = -11.0768
$ws.Range("C70").Value = -12.6775
$ws.Range("D70").Value = -8.002100000000002
$ws.Range("A76").Value = -19.27469999999999
$ws.Range("C76").Value = -13.04640000000001
$ws.Range("A78").Value = -19.77989999999998
$ws.Range("C81").Value = -13.7684
$ws.Range("D81").Value = -8.143300000000005
$ws.Range("A82").Value = -22.287
$ws.Range("A83").Value = -22.00809999999999
$ws.Range("C83").Value = -13.9425
$ws.Range("E84").Value = 16.95509999999999
$ws.Range("D88").Value = -7.272899999999993
$ws.Range("A93").Value = -20.41389999999998
$ws.Range("C94").Value = -10.3852
$ws.Range("D96").Value = -8.219399999999993
$ws.Range("A97").Value = -21.9288
$ws.Range("D99").Value = -7.611499999999991
$ws.Range("C100").Value = -10.83310000000001
$ws.Range("E100").Value = 16.65349999999999
$ws.Range("E101").Value = 16.82020000000001
$ws.Range("C102").Value = -13.37400000000001